$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the exact formatting of the existing
# header cells (e.g. G1) via copy/paste-format so no new style entry
# is created.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add value 0 in H2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 0
